$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Add the three new header cells in row 1 (new language / ordering columns)
    # Order of first-use controls the shared string table index, so write
    # U1, then W1, then V1 to reproduce indices 106/107/108 respectively.
    $ws.Range("U1").Value2 = "drseca_najprej"
    $ws.Range("W1").Value2 = "leva_os_en"
    $ws.Range("V1").Value2 = "leva_os_si"

    # Update the selection shown for this sheet
    $ws.Range("U1:W1").Select() | Out-Null
}

# Make Sheet1 the active sheet / tab (selected last so it sticks)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("U1:W1").Select() | Out-Null
